$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns for data rows stay text,
# matching the source formatting (e.g. "135.90", "1.015", "30.702.69").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.702.69"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "2.122.32"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +1.35%  "

$ws.Range("D5").Value = "338.89"
$ws.Range("E5").Value = "  +1.95%  "

$ws.Range("D7").Value = "0.5284"
$ws.Range("E7").Value = "  +1.06%  "

$ws.Range("D8").Value = "0.4555"
$ws.Range("E8").Value = "  +1.53%  "

$ws.Range("D9").Value = "54.19"
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("D10").Value = "0.09123"
$ws.Range("E10").Value = "  +1.92%  "

$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("D12").Value = "24.51"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").Value = "2.123.72"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").Value = "6.857"
$ws.Range("E14").Value = "  +1.44%  "

$ws.Range("D15").Value = "8.114"
$ws.Range("E15").Value = "  +4.51%  "

$ws.Range("D16").Value = "98.87"
$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("D17").Value = "0.00001172"
$ws.Range("E17").Value = "  +4.15%  "

$ws.Range("D18").Value = "1.015"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("D19").Value = "0.06712"
$ws.Range("E19").Value = "  +1.60%  "

$ws.Range("D20").Value = "19.57"
$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("D21").Value = "1.013"
$ws.Range("E21").Value = "  +1.23%  "

$ws.Range("D22").Value = "6.465"
$ws.Range("E22").Value = "  +2.65%  "

$ws.Range("D23").Value = "30.778.08"
$ws.Range("E23").Value = "  +0.81%  "

$ws.Range("D24").Value = "12.97"
$ws.Range("E24").Value = "  +5.22%  "

$ws.Range("D25").Value = "2.378"
$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("D26").Value = "2.371.34"
$ws.Range("E26").Value = "  +1.13%  "

$ws.Range("D27").Value = "22.53"
$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("D28").Value = "165.99"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("D29").Value = "2.557"
$ws.Range("E29").Value = "  -1.27%  "

$ws.Range("D30").Value = "135.90"
$ws.Range("E30").Value = "  +2.36%  "

$ws.Range("D31").Value = "1.212"
$ws.Range("E31").Value = "  +0.82%  "

$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").Value = "6.413"
$ws.Range("E33").Value = "  +4.06%  "

$ws.Range("D34").Value = "1.643"
$ws.Range("E34").Value = "  -2.14%  "

$ws.Range("D35").Value = "3.957"
$ws.Range("E35").Value = "  +0.52%  "

$ws.Range("D36").Value = "10.59"
$ws.Range("E36").Value = "  +1.75%  "

$ws.Range("D37").Value = "6.001"
$ws.Range("E37").Value = "  +9.24%  "

$ws.Range("D38").Value = "0.02666"
$ws.Range("E38").Value = "  +3.68%  "

$ws.Range("D39").Value = "0.06890"
$ws.Range("E39").Value = "  +1.53%  "

$ws.Range("E40").Value = "  +2.01%  "

$ws.Range("D41").Value = "12.62"
$ws.Range("E41").Value = "  -1.95%  "

$ws.Range("D42").Value = "0.6918"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").Value = "1.267"
$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("D44").Value = "15.20"
$ws.Range("E44").Value = "  +9.22%  "

$ws.Range("D45").Value = "0.6485"
$ws.Range("E45").Value = "  +1.73%  "

$ws.Range("D46").Value = "2.315"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("D47").Value = "0.00000000370"
$ws.Range("E47").Value = "  +15.95%  "

$ws.Range("D48").Value = "3.707"
$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("D49").Value = "1.261"
$ws.Range("E49").Value = "  +1.20%  "

$ws.Range("D50").Value = "83.26"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").Value = "0.07318"
$ws.Range("E51").Value = "  +3.78%  "

